$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (copying cell formatting from the column to its left),
# shifting the old "dbExcel"/"WebExcel" columns (B,C) one position to the right (C,D).
[void]$ws.Columns.Item(2).Insert(0)

# Populate the new StatQuery header + query cells.
$ws.Range("B1").Value2 = "StatQuery"
$ws.Range("B2").Value2 = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (single lung lobe)']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the wrapped-text styling used by column A's query cell.
$ws.Range("B2").WrapText = $true

# Keep the new column the same width as column A.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Update the active selection to reflect the edited cell.
[void]$ws.Range("B2").Select()
